$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '72.887.07'
$ws.Range("E2").Value = '  +1.02%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '4.013.07'
$ws.Range("E3").Value = '  -0.27%  '
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '592.45'
$ws.Range("E5").Value = '  +12.29%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '153.02'
$ws.Range("E6").Value = '  +1.48%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.688'
$ws.Range("E7").Value = '  -1.68%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '1.00'
$ws.Range("E8").Value = '  -0.04%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.759'
$ws.Range("E9").Value = '  +1.37%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.170'
$ws.Range("E10").Value = '  -0.87%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '54.20'
$ws.Range("E11").Value = '  +8.79%  '
$ws.Range("E12").Value = '  -1.76%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '10.99'
$ws.Range("E13").Value = '  +2.99%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '4.660.12'
$ws.Range("E14").Value = '  -0.17%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.022.33'
$ws.Range("E15").Value = '  -0.19%  '
$ws.Range("E16").Value = '  +7.66%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '14.24'
$ws.Range("E17").Value = '  +0.99%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '20.67'
$ws.Range("E18").Value = '  +0.35%  '
$ws.Range("E19").Value = '  -0.50%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '72.751.76'
$ws.Range("E20").Value = '  +0.83%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '438.32'
$ws.Range("E21").Value = '  +1.36%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.82'
$ws.Range("E22").Value = '  +15.14%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '97.06'
$ws.Range("E23").Value = '  -1.04%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.52'
$ws.Range("E24").Value = '  +0.78%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '14.40'
$ws.Range("E25").Value = '  +1.01%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '4.39'
$ws.Range("E26").Value = '  +19.33%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.42'
$ws.Range("E27").Value = '  +1.32%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.78'
$ws.Range("E28").Value = '  +0.70%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.92'
$ws.Range("E29").Value = '  +1.11%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '36.72'
$ws.Range("E30").Value = '  +0.00%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '7.92'
$ws.Range("E31").Value = '  +7.35%  '
$ws.Range("B32").Value = 'Cosmos'
$ws.Range("C32").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '13.72'
$ws.Range("E32").Value = '  +1.90%  '
$ws.Range("B33").Value = 'Hedera'
$ws.Range("C33").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.133'
$ws.Range("E33").Value = '  +1.26%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '49.89'
$ws.Range("E34").Value = '  +3.96%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '685.34'
$ws.Range("E35").Value = '  +0.75%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '71.88'
$ws.Range("E36").Value = '  +9.52%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.443'
$ws.Range("E37").Value = '  -1.31%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0₃0870'
$ws.Range("E38").Value = '  +5.42%  '
$ws.Range("E39").Value = '  -1.28%  '
$ws.Range("E40").Value = '  +3.27%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '11.19'
$ws.Range("E41").Value = '  +9.13%  '
$ws.Range("B42").Value = 'Dai'
$ws.Range("C42").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.00'
$ws.Range("E42").Value = '  -0.06%  '
$ws.Range("B43").Value = 'ThetaToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.34'
$ws.Range("E43").Value = '  -1.31%  '
$ws.Range("B44").Value = 'FirstDigitalUSD'
$ws.Range("C44").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '1.00'
$ws.Range("E44").Value = '  +0.03%  '
$ws.Range("B45").Value = 'VeChain'
$ws.Range("C45").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0493'
$ws.Range("E45").Value = '  -0.14%  '
$ws.Range("E46").Value = '  +2.98%  '
$ws.Range("E47").Value = '  +0.45%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.36'
$ws.Range("E48").Value = '  -0.13%  '
$ws.Range("E49").Value = '  +8.09%  '
$ws.Range("E50").Value = '  +0.22%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.843.55'
$ws.Range("E51").Value = '  +11.35%  '
